$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update row 12 (task: "Escribir el código para que el usuario pueda introducir
# por teclado el nombre de usuario y la contraseña") with new ENCARGADO,
# TIEMPO ESTIMADO and TIEMPO INVERTIDO values.
$ws.Range("D12").Value = "Enrique, Sergio"
$ws.Range("E12").Value = "1,5 horas"
$ws.Range("F12").Value = "2 horas"

# Update the view/selection state to match the saved workbook
# (topLeftCell -> C2, activeCell -> F13).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F13").Select()
